# Commit message: "Fruta / hortaliza, semanal"
# The diff shows a new weekly price-report row being inserted at row 61
# of the (single) worksheet; every existing row from 61 downwards shifts
# down by one, and the previously-last row (old row 136) becomes the new
# last row (137). The sheet's dimension grows from A1:R136 to A1:R137.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 61 - this shifts rows 61..136 down to 62..137
# and preserves existing cell formatting (e.g. the date style on column D).
$ws.Rows(61).Insert()

# Populate the newly inserted row 61 with the new weekly record.
$ws.Range("A61").Value = 3
$ws.Range("B61").Value = "Femacal de La Calera"
$ws.Range("C61").Value = "Coquimbo"
$ws.Range("D61").Value = 44601
$ws.Range("E61").Value = 5
$ws.Range("F61").Value = 100112052
$ws.Range("G61").Value = "Albahaca"
$ws.Range("H61").Value = "Sin especificar"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 60
$ws.Range("K61").Value = 4000
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = 4000
$ws.Range("N61").Value = "$/docena de matas"
$ws.Range("O61").Value = "Provincia de Quillota"
$ws.Range("P61").Value = 667
$ws.Range("Q61").Value = 6
$ws.Range("R61").Value = "Hortaliza"
